$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the three new cells: B1, A2 (numeric 0, styled) and B2 (text label)
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the "bold, centered, top-aligned, thin-bordered" style on B1 ...
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4160     # xlTop
$r1.Borders.LineStyle = 1         # xlContinuous
$r1.Borders.Weight = 2            # xlThin

# ... and copy the same formatting onto A2 so both cells end up sharing
# the exact same style definition (avoids creating a duplicate style entry).
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
